# Apply the 2025-11-12 Betfair back/lay refresh:
#  - insert a new 'Friendly Matches' fixture as row 2 (pushes the four existing
#    fixtures down from rows 2-5 to rows 3-6)
#  - refresh the odds figures for the fixtures that moved to rows 3-6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2-5) down to (3-6) by inserting a blank row at 2.
# ClearFormats() drops the formatting Excel copies down from the header row on
# insert, so the new row matches the plain (unstyled) look of the other data rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Column B (Date) is the literal text '2025-11-12' on every data row. Copy the
# already-shifted value from B3 into B2 instead of retyping the string, otherwise
# Excel's COM layer auto-parses the ISO-looking text into a date serial number.
$ws.Range("B3").Copy($ws.Range("B2"))

# Row 2: Friendly Matches | SV Lafnitz vs KSV 1919
$ws.Range("A2").Value = 'Friendly Matches'
$ws.Range("C2").Value = '13:00:00'
$ws.Range("D2").Value = 'SV Lafnitz'
$ws.Range("E2").Value = 'KSV 1919'
$ws.Range("F2").Value = 1.04
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.04
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.02
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.01
$ws.Range("O2").Value = 1.09
$ws.Range("P2").Value = 1.08
$ws.Range("Q2").Value = 1.09
$ws.Range("R2").Value = 1.08
$ws.Range("S2").Value = 1.09
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3: Bosnian Premier League | Borac Banja Luka vs Zrinjski
$ws.Range("A3").Value = 'Bosnian Premier League'
$ws.Range("C3").Value = '14:00:00'
$ws.Range("D3").Value = 'Borac Banja Luka'
$ws.Range("E3").Value = 'Zrinjski'
$ws.Range("F3").Value = 1.85
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 1.09
$ws.Range("I3").Value = 5.5
$ws.Range("J3").Value = 1.2
$ws.Range("K3").Value = 7.6
$ws.Range("L3").Value = 1.27
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.34
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.34
$ws.Range("Q3").Value = 1.42
$ws.Range("R3").Value = 1.2
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04
$ws.Range("V3").Value = 1.22
$ws.Range("W3").Value = 1.5
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4: Colombian Primera B | Real Cartagena vs Real Soacha Cundinamarca FC
$ws.Range("A4").Value = 'Colombian Primera B'
$ws.Range("C4").Value = '20:10:00'
$ws.Range("D4").Value = 'Real Cartagena'
$ws.Range("E4").Value = 'Real Soacha Cundinamarca FC'
$ws.Range("F4").Value = 1.87
$ws.Range("G4").Value = 2.26
$ws.Range("H4").Value = 3.9
$ws.Range("I4").Value = 6.2
$ws.Range("J4").Value = 2.9
$ws.Range("K4").Value = 4
$ws.Range("L4").Value = 1.46
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 2.68
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 1.67
$ws.Range("Q4").Value = 2.02
$ws.Range("R4").Value = 1.25
$ws.Range("S4").Value = 3.45
$ws.Range("T4").Value = 1.94
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.19
$ws.Range("W4").Value = 1.79
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 5: Brazilian Serie A | Atletico MG vs Fortaleza EC
$ws.Range("A5").Value = 'Brazilian Serie A'
$ws.Range("C5").Value = '20:30:00'
$ws.Range("D5").Value = 'Atletico MG'
$ws.Range("E5").Value = 'Fortaleza EC'
$ws.Range("F5").Value = 1.74
$ws.Range("G5").Value = 1.75
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 6.2
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 3.9
$ws.Range("L5").Value = 1.45
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 3.45
$ws.Range("O5").Value = 1.37
$ws.Range("P5").Value = 1.82
$ws.Range("Q5").Value = 2.12
$ws.Range("R5").Value = 1.31
$ws.Range("S5").Value = 3.95
$ws.Range("T5").Value = 2.08
$ws.Range("U5").Value = 1.87
$ws.Range("V5").Value = 1.19
$ws.Range("W5").Value = 2.34
$ws.Range("X5").Value = 13
$ws.Range("Y5").Value = 18
$ws.Range("Z5").Value = 46
$ws.Range("AA5").Value = 170
$ws.Range("AB5").Value = 7.6
$ws.Range("AC5").Value = 8.4
$ws.Range("AD5").Value = 23
$ws.Range("AE5").Value = 95
$ws.Range("AF5").Value = 9.199999999999999
$ws.Range("AG5").Value = 10.5
$ws.Range("AH5").Value = 24
$ws.Range("AI5").Value = 100
$ws.Range("AJ5").Value = 17
$ws.Range("AK5").Value = 20
$ws.Range("AL5").Value = 42
$ws.Range("AM5").Value = 150
$ws.Range("AN5").Value = 12.5
$ws.Range("AO5").Value = 130

# Row 6: Colombian Primera A | Boyaca Chico vs Millonarios
$ws.Range("A6").Value = 'Colombian Primera A'
$ws.Range("C6").Value = '22:20:00'
$ws.Range("D6").Value = 'Boyaca Chico'
$ws.Range("E6").Value = 'Millonarios'
$ws.Range("F6").Value = 5.3
$ws.Range("G6").Value = 6.6
$ws.Range("H6").Value = 1.75
$ws.Range("I6").Value = 1.89
$ws.Range("J6").Value = 3.35
$ws.Range("K6").Value = 4.1
$ws.Range("L6").Value = 1.49
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 2.62
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 1.64
$ws.Range("Q6").Value = 2.24
$ws.Range("R6").Value = 1.23
$ws.Range("S6").Value = 3.9
$ws.Range("T6").Value = 2.06
$ws.Range("U6").Value = 1.74
$ws.Range("V6").Value = 2.12
$ws.Range("W6").Value = 1.19
$ws.Range("X6").Value = 980
$ws.Range("Y6").Value = 980
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 980
$ws.Range("AB6").Value = 980
$ws.Range("AC6").Value = 980
$ws.Range("AD6").Value = 980
$ws.Range("AE6").Value = 980
$ws.Range("AF6").Value = 980
$ws.Range("AG6").Value = 980
$ws.Range("AH6").Value = 980
$ws.Range("AI6").Value = 65
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 130
$ws.Range("AL6").Value = 140
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 980
